$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 175. This shifts the
# existing rows 175:266 down to 176:267 and extends the used range to
# A1:R267, matching the target dimension.
$ws.Rows("175:175").Insert()

# Populate the newly inserted row 175 with the new weekly price record.
$ws.Range("A175").Value = 5
$ws.Range("B175").Value = "Macroferia Regional de Talca"
$ws.Range("C175").Value = "Maule"
$ws.Range("D175").Value = 44572
$ws.Range("E175").Value = 7
$ws.Range("F175").Value = 100112032
$ws.Range("G175").Value = "Zapallo italiano"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 200
$ws.Range("K175").Value = 7000
$ws.Range("L175").Value = 7000
$ws.Range("M175").Value = 7000
$ws.Range("N175").Value = "$/caja 50 unidades"
$ws.Range("O175").Value = "Región del Maule"
$ws.Range("P175").Value = 140
$ws.Range("Q175").Value = 50
$ws.Range("R175").Value = "Hortaliza"
